$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1949.5
$ws.Range("I2").Value = 899
$ws.Range("K2").Value = 899
$ws.Range("M2").Value = -786
$ws.Range("H9").Value = 8875.154
$ws.Range("J9").Value = 116
$ws.Range("L9").Value = 116
$ws.Range("N9").Value = -454
$ws.Range("H21").Value = 17503.75
$ws.Range("J21").Value = 19999.334
$ws.Range("L21").Value = 19999.334
$ws.Range("N21").Value = -20935.334
$ws.Range("H23").Value = 17503.75
$ws.Range("J23").Value = 19999.334
$ws.Range("L23").Value = 19999.334
$ws.Range("N23").Value = -20467.334
$ws.Range("H74").Value = 10132.267
$ws.Range("I74").Value = 9070.286
$ws.Range("K74").Value = 9070.286
$ws.Range("M74").Value = -8134.286
$ws.Range("H77").Value = 10132.267
$ws.Range("I77").Value = 9070.286
$ws.Range("K77").Value = 45351.43
$ws.Range("M77").Value = -40671.43
$ws.Range("H86").Value = 2197886.8
$ws.Range("J86").Value = 3763726.8
$ws.Range("L86").Value = 3763726.8
$ws.Range("N86").Value = -3765972.8
$ws.Range("H89").Value = 2197886.8
$ws.Range("J89").Value = 3763726.8
$ws.Range("L89").Value = 18818634
$ws.Range("N89").Value = -18829866
$ws.Range("H135").Value = 4682.207
$ws.Range("J135").Value = 9166.333000000001
$ws.Range("L135").Value = 82496.997
$ws.Range("N135").Value = -87566.997
$ws.Range("H138").Value = 5545.6514
$ws.Range("I138").Value = 1439.579
$ws.Range("J138").Value = 8796.291999999999
$ws.Range("K138").Value = 4318.737
$ws.Range("L138").Value = 26388.876
$ws.Range("M138").Value = 821.2629999999999
$ws.Range("N138").Value = -36668.876

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2735.486
$ws.Range("I32").Value = 2684.5796
$ws.Range("K32").Value = 2684.5796
$ws.Range("M32").Value = -2397.5796
$ws.Range("H101").Value = 56665
$ws.Range("J101").Value = 56665
$ws.Range("L101").Value = 56665
$ws.Range("N101").Value = -63155
$ws.Range("H110").Value = 459897.28
$ws.Range("I110").Value = 629107.75
$ws.Range("K110").Value = 629107.75
$ws.Range("M110").Value = -627062.75
$ws.Range("H132").Value = 247655.8
$ws.Range("I132").Value = 300973
$ws.Range("J132").Value = 7728.375
$ws.Range("K132").Value = 902919
$ws.Range("L132").Value = 23185.125
$ws.Range("M132").Value = -900389
$ws.Range("N132").Value = -28245.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1921.5
$ws.Range("J20").Value = 2647.6667
$ws.Range("L20").Value = 2647.6667
$ws.Range("N20").Value = -3141.6667
$ws.Range("H105").Value = 7578632
$ws.Range("I105").Value = 2202.5
$ws.Range("K105").Value = 2202.5
$ws.Range("M105").Value = -455.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1337770.9
$ws.Range("I31").Value = 16949152
$ws.Range("J31").Value = 36822.5
$ws.Range("K31").Value = 16949152
$ws.Range("L31").Value = 36822.5
$ws.Range("M31").Value = -16948857
$ws.Range("N31").Value = -37412.5
$ws.Range("H34").Value = 1337770.9
$ws.Range("I34").Value = 16949152
$ws.Range("J34").Value = 36822.5
$ws.Range("K34").Value = 16949152
$ws.Range("L34").Value = 36822.5
$ws.Range("M34").Value = -16948950
$ws.Range("N34").Value = -37226.5
$ws.Range("H62").Value = 3600.2
$ws.Range("J62").Value = 4500
$ws.Range("L62").Value = 4500
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 3600.2
$ws.Range("J65").Value = 4500
$ws.Range("L65").Value = 22500
$ws.Range("N65").Value = -28740
$ws.Range("H134").Value = 232467.27
$ws.Range("I134").Value = 2457.2942
$ws.Range("J134").Value = 1014501.2
$ws.Range("K134").Value = 7371.882599999999
$ws.Range("L134").Value = 3043503.6
$ws.Range("M134").Value = -4836.882599999999
$ws.Range("N134").Value = -3048573.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 874.6
$ws.Range("J75").Value = 357.66666
$ws.Range("L75").Value = 1072.99998
$ws.Range("N75").Value = -3068.99998
$ws.Range("H78").Value = 874.6
$ws.Range("J78").Value = 357.66666
$ws.Range("L78").Value = 3218.99994
$ws.Range("N78").Value = -13202.99994
$ws.Range("H87").Value = 5666.3335
$ws.Range("I87").Value = 5666.3335
$ws.Range("K87").Value = 16999.0005
$ws.Range("M87").Value = -15751.0005
$ws.Range("H90").Value = 5666.3335
$ws.Range("I90").Value = 5666.3335
$ws.Range("K90").Value = 50997.0015
$ws.Range("M90").Value = -44757.0015
$ws.Range("H122").Value = 25393456
$ws.Range("I122").Value = 59218960
$ws.Range("J122").Value = 24325.75
$ws.Range("K122").Value = 532970640
$ws.Range("L122").Value = 218931.75
$ws.Range("M122").Value = -532968190
$ws.Range("N122").Value = -223831.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 1853335.1
$ws.Range("J10").Value = 6669.3335
$ws.Range("L10").Value = 6669.3335
$ws.Range("N10").Value = -7007.3335
$ws.Range("H12").Value = 10004
$ws.Range("J12").Value = 10004
$ws.Range("L12").Value = 10004
$ws.Range("N12").Value = -10284
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H132").Value = 16680.646
$ws.Range("I132").Value = 1879.3334
$ws.Range("K132").Value = 5638.0002
$ws.Range("M132").Value = -3108.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 4034.5
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 4821.4
$ws.Range("K11").Value = 100
$ws.Range("L11").Value = 4821.4
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = -5101.4
$ws.Range("H22").Value = 2496.3845
$ws.Range("I22").Value = 2781.3333
$ws.Range("J22").Value = 2252.1428
$ws.Range("K22").Value = 2781.3333
$ws.Range("L22").Value = 2252.1428
$ws.Range("M22").Value = -2486.3333
$ws.Range("N22").Value = -2842.1428
$ws.Range("H25").Value = 10526
$ws.Range("J25").Value = 10526
$ws.Range("L25").Value = 10526
$ws.Range("N25").Value = -10986
$ws.Range("H27").Value = 2496.3845
$ws.Range("I27").Value = 2781.3333
$ws.Range("J27").Value = 2252.1428
$ws.Range("K27").Value = 2781.3333
$ws.Range("L27").Value = 2252.1428
$ws.Range("M27").Value = -2674.3333
$ws.Range("N27").Value = -2466.1428
$ws.Range("H100").Value = 120793.89
$ws.Range("I100").Value = 169441
$ws.Range("K100").Value = 169441
$ws.Range("M100").Value = -168900
$ws.Range("H106").Value = 43369.668
$ws.Range("J106").Value = 43369.668
$ws.Range("L106").Value = 43369.668
$ws.Range("N106").Value = -45893.668
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 4999.5
$ws.Range("I24").Value = 4999.5
$ws.Range("K24").Value = 4999.5
$ws.Range("M24").Value = -4769.5
$ws.Range("H88").Value = 10000
$ws.Range("J88").Value = 10000
$ws.Range("L88").Value = 10000
$ws.Range("N88").Value = -10812
$ws.Range("H91").Value = 10000
$ws.Range("J91").Value = 10000
$ws.Range("L91").Value = 10000
$ws.Range("N91").Value = -12808
$ws.Range("H104").Value = 66624.75
$ws.Range("J104").Value = 66624.75
$ws.Range("L104").Value = 66624.75
$ws.Range("N104").Value = -73612.75
$ws.Range("H122").Value = 4488.381
$ws.Range("I122").Value = 3250.5386
$ws.Range("K122").Value = 9751.6158
$ws.Range("M122").Value = -7301.6158
